# Reproduce a blank-column issue (see issue #35) by adding a new test case
# to the fancy-test workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet2: add a couple of helper/label cells in column E, and widen
#     column D so the issue (column D looking blank / too narrow) is visible.
$ws2 = $wb.Worksheets.Item("Sheet2")

# ~2x the default column width, so column D (which looked blank/too-narrow
# per issue #35) is clearly wide enough.
$ws2.Columns.Item(4).ColumnWidth = 16

$ws2.Range("E2").Value = "< that is column D"
$ws2.Range("E3").Value = "< it should be ~2x wide"
$ws2.Range("E4").Value = "< as column E"

[void]$ws2.Range("E4").Select()

# --- Add a new Sheet3 (after the last existing sheet) with a small repro
#     case for the issue.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

$ws3.Range("B1").Value = "test"
$ws3.Range("B2").Value = "test"
$ws3.Range("B3").Value = 123

[void]$ws3.Range("B4").Select()
[void]$ws3.Activate()
